$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 107 (shifts all existing rows 107-229 down to 108-230,
# preserving all of their data/styles), then populate the new row with the new weekly
# price-report record.
$ws.Rows.Item(107).Insert()

$ws.Range("A107").Value = 4
$ws.Range("B107").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C107").Value = "Los Lagos"
$ws.Range("D107").Value = 44740
$ws.Range("E107").Value = 10
$ws.Range("F107").Value = 100112039
$ws.Range("G107").Value = "Ciboulette"
$ws.Range("H107").Value = "Sin especificar"
$ws.Range("I107").Value = "Primera"
$ws.Range("J107").Value = 240
$ws.Range("K107").Value = 2500
$ws.Range("L107").Value = 2500
$ws.Range("M107").Value = 2500
$ws.Range("N107").Value = "$/docena de atados"
$ws.Range("O107").Value = "Región Metropolitana"
$ws.Range("P107").Value = 833
$ws.Range("Q107").Value = 3
$ws.Range("R107").Value = "Hortaliza"
